# Update CDA Logical model for ST.r2b
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInclude = $wb.Worksheets.Item("Include from ActClass")

# Rename the "Include from ActClass" sheet to "Include #0"
$wsInclude.Name = "Include #0"

# Insert a new row for "Jurisdiction" after the "Contact" row (row 10),
# pushing "Description", "Purpose", "Copyright", "Immutable" down by one.
$wsMeta.Rows.Item(11).Insert()

# Copy the style of the row above (Contact row) onto the freshly inserted row
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsMeta.Range("A11").Value = "Jurisdiction"
# B11 is left blank (no value), matching the "Jurisdiction" row's empty
# Value column.

# Update Version value
$wsMeta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value
$wsMeta.Range("B8").Value = "2025-10-29T22:15:57+01:00"
